$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" '63.212.53'
Set-TextCell "E2" '  +0.09%  '

# Row 3
Set-TextCell "D3" '3.244.29'
Set-TextCell "E3" '  +0.49%  '

# Row 4
Set-TextCell "D4" '0.999'
Set-TextCell "E4" '  -0.31%  '

# Row 5
Set-TextCell "D5" '530.50'
Set-TextCell "E5" '  +4.70%  '

# Row 6
Set-TextCell "D6" '171.75'
Set-TextCell "E6" '  -1.12%  '

# Row 7
Set-TextCell "D7" '0.597'
Set-TextCell "E7" '  +2.77%  '

# Row 8
Set-TextCell "E8" '  -0.03%  '

# Row 9
Set-TextCell "D9" '3.243.52'
Set-TextCell "E9" '  +0.79%  '

# Row 10
Set-TextCell "D10" '0.607'
Set-TextCell "E10" '  +0.47%  '

# Row 11
Set-TextCell "D11" '53.21'
Set-TextCell "E11" '  -5.30%  '

# Row 12
Set-TextCell "E12" '  +5.09%  '

# Row 13
Set-TextCell "D13" '0.0000255'
Set-TextCell "E13" '  +2.39%  '

# Row 14
Set-TextCell "D14" '9.15'
Set-TextCell "E14" '  +2.77%  '

# Row 15
Set-TextCell "D15" '3.759.15'
Set-TextCell "E15" '  +0.22%  '

# Row 16
Set-TextCell "E16" '  -0.45%  '

# Row 17
Set-TextCell "D17" '3.226.86'
Set-TextCell "E17" '  -0.15%  '

# Row 18
Set-TextCell "B18" 'Chainlink'
Set-TextCell "C18" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell "D18" '17.26'
Set-TextCell "E18" '  +2.41%  '

# Row 19
Set-TextCell "B19" 'WrappedBTC'
Set-TextCell "C19" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell "D19" '63.022.28'
Set-TextCell "E19" '  +0.17%  '

# Row 20
Set-TextCell "D20" '11.12'
Set-TextCell "E20" '  +5.33%  '

# Row 21
Set-TextCell "D21" '0.970'
Set-TextCell "E21" '  +4.97%  '

# Row 22
Set-TextCell "D22" '367.03'
Set-TextCell "E22" '  +1.28%  '

# Row 23
Set-TextCell "D23" '3.76'
Set-TextCell "E23" '  +5.86%  '

# Row 24
Set-TextCell "D24" '81.08'
Set-TextCell "E24" '  +3.36%  '

# Row 25
Set-TextCell "D25" '11.20'
Set-TextCell "E25" '  +4.94%  '

# Row 26
Set-TextCell "D26" '3.99'
Set-TextCell "E26" '  +7.75%  '

# Row 27
Set-TextCell "E27" '  -0.42%  '

# Row 28
Set-TextCell "D28" '2.65'
Set-TextCell "E28" '  +2.39%  '

# Row 29
Set-TextCell "D29" '11.31'
Set-TextCell "E29" '  +2.83%  '

# Row 30
Set-TextCell "D30" '8.22'
Set-TextCell "E30" '  +1.34%  '

# Row 31
Set-TextCell "D31" '28.53'
Set-TextCell "E31" '  +2.52%  '

# Row 32
Set-TextCell "D32" '634.82'
Set-TextCell "E32" '  +0.55%  '

# Row 33
Set-TextCell "D33" '6.48'
Set-TextCell "E33" '  -0.37%  '

# Row 34
Set-TextCell "D34" '11.22'
Set-TextCell "E34" '  +3.31%  '

# Row 35
Set-TextCell "D35" '0.106'
Set-TextCell "E35" '  +4.92%  '

# Row 36
Set-TextCell "D36" '56.90'
Set-TextCell "E36" '  -2.55%  '

# Row 37
Set-TextCell "E37" '  +0.17%  '

# Row 38
Set-TextCell "D38" '36.80'
Set-TextCell "E38" '  +5.07%  '

# Row 39
Set-TextCell "D39" '0.378'
Set-TextCell "E39" '  +2.52%  '

# Row 40
Set-TextCell "B40" 'PEPE'
Set-TextCell "C40" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell "D40" '0.0₃0722'
Set-TextCell "E40" '  +14.60%  '

# Row 41
Set-TextCell "B41" 'FirstDigitalUSD'
Set-TextCell "C41" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell "D41" '0.996'
Set-TextCell "E41" '  -0.21%  '

# Row 42
Set-TextCell "D42" '0.123'
Set-TextCell "E42" '  +3.10%  '

# Row 43
Set-TextCell "E43" '  +12.31%  '

# Row 44
Set-TextCell "D44" '2.880.17'
Set-TextCell "E44" '  +2.02%  '

# Row 45
Set-TextCell "E45" '  +6.37%  '

# Row 46
Set-TextCell "D46" '2.68'
Set-TextCell "E46" '  +4.50%  '

# Row 47
Set-TextCell "B47" 'VeChain'
Set-TextCell "C47" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell "D47" '0.0395'
Set-TextCell "E47" '  +5.85%  '

# Row 48
Set-TextCell "B48" 'ApeXProtocol'
Set-TextCell "C48" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell "D48" '3.08'
Set-TextCell "E48" '  +6.75%  '

# Row 49
Set-TextCell "D49" '2.59'
Set-TextCell "E49" '  -0.27%  '

# Row 50
Set-TextCell "E50" '  +3.50%  '

# Row 51
Set-TextCell "D51" '133.95'
Set-TextCell "E51" '  +2.30%  '
